# Updated cryptos list with GitHub Actions: refresh Price/Volume(1h) figures
# and roll the ranked coin list down one slot (a new coin, OKB, enters at
# rank 9 / row 10, pushing RocketPoolETH off the bottom of the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '36.567.29'
$ws.Cells.Item(2, 5).Value = '  +1.06%  '
$ws.Cells.Item(3, 4).Value = '1.956.44'
$ws.Cells.Item(3, 5).Value = '  -0.06%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).Formula = "'244.56"
$ws.Cells.Item(5, 5).Value = '  +1.09%  '
$ws.Cells.Item(6, 4).Formula = "'0.621"
$ws.Cells.Item(6, 5).Value = '  +0.55%  '
$ws.Cells.Item(7, 4).Formula = "'58.67"
$ws.Cells.Item(7, 5).Value = '  +1.78%  '
$ws.Cells.Item(8, 5).Value = '  -0.09%  '
$ws.Cells.Item(9, 5).Value = '  +0.10%  '
$ws.Cells.Item(10, 2).Value = 'OKB'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(10, 4).Formula = "'56.30"
$ws.Cells.Item(10, 5).Value = '  -0.78%  '
$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).Formula = "'0.0862"
$ws.Cells.Item(11, 5).Value = '  +10.10%  '
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).Formula = "'0.104"
$ws.Cells.Item(12, 5).Value = '  +1.48%  '
$ws.Cells.Item(13, 2).Value = 'Avalanche'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(13, 4).Formula = "'22.54"
$ws.Cells.Item(13, 5).Value = '  +4.57%  '
$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(14, 4).Formula = "'0.832"
$ws.Cells.Item(14, 5).Value = '  -1.04%  '
$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).Value = '2.242.34'
$ws.Cells.Item(15, 5).Value = '  -0.14%  '
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).Formula = "'13.72"
$ws.Cells.Item(16, 5).Value = '  -0.93%  '
$ws.Cells.Item(17, 2).Value = 'Polkadot'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(17, 4).Formula = "'5.26"
$ws.Cells.Item(17, 5).Value = '  -1.80%  '
$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value = '1.962.12'
$ws.Cells.Item(18, 5).Value = '  +0.30%  '
$ws.Cells.Item(19, 2).Value = 'WrappedBTC'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(19, 4).Value = '36.484.98'
$ws.Cells.Item(19, 5).Value = '  +1.21%  '
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(20, 4).Value = '0.0₃0882'
$ws.Cells.Item(20, 5).Value = '  +4.43%  '
$ws.Cells.Item(21, 2).Value = 'Litecoin'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(21, 4).Formula = "'70.14"
$ws.Cells.Item(21, 5).Value = '  -0.94%  '
$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Formula = "'230.17"
$ws.Cells.Item(22, 5).Value = '  -2.45%  '
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(23, 4).Formula = "'5.10"
$ws.Cells.Item(23, 5).Value = '  -1.17%  '
$ws.Cells.Item(24, 2).Value = 'Dai'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(24, 4).Formula = "'1.00"
$ws.Cells.Item(24, 5).Value = '  +0.08%  '
$ws.Cells.Item(25, 2).Value = 'PancakeSwap'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(25, 4).Formula = "'2.49"
$ws.Cells.Item(25, 5).Value = '  -1.05%  '
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(26, 4).Formula = "'2.31"
$ws.Cells.Item(26, 5).Value = '  +1.98%  '
$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(27, 4).Formula = "'9.44"
$ws.Cells.Item(27, 5).Value = '  -1.67%  '
$ws.Cells.Item(28, 2).Value = 'Monero'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(28, 4).Formula = "'162.62"
$ws.Cells.Item(28, 5).Value = '  +1.38%  '
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(29, 4).Formula = "'0.136"
$ws.Cells.Item(29, 5).Value = '  +11.09%  '
$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).Formula = "'19.66"
$ws.Cells.Item(30, 5).Value = '  +0.00%  '
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).Formula = "'0.118"
$ws.Cells.Item(31, 5).Value = '  +0.06%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).Formula = "'1.19"
$ws.Cells.Item(32, 5).Value = '  +6.65%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).Formula = "'4.74"
$ws.Cells.Item(33, 5).Value = '  -1.70%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Formula = "'0.0643"
$ws.Cells.Item(34, 5).Value = '  +5.82%  '
$ws.Cells.Item(35, 5).Value = '  -1.15%  '
$ws.Cells.Item(36, 2).Value = 'THORChain'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(36, 4).Formula = "'6.43"
$ws.Cells.Item(36, 5).Value = '  +7.95%  '
$ws.Cells.Item(37, 2).Value = 'BinanceUSD'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(37, 4).Formula = "'1.00"
$ws.Cells.Item(37, 5).Value = '  +0.11%  '
$ws.Cells.Item(38, 5).Value = '  -2.02%  '
$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(39, 4).Formula = "'2.20"
$ws.Cells.Item(39, 5).Value = '  -2.47%  '
$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(40, 4).Formula = "'3.05"
$ws.Cells.Item(40, 5).Value = '  +2.80%  '
$ws.Cells.Item(41, 2).Value = 'Cronos'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(41, 4).Formula = "'0.1000"
$ws.Cells.Item(41, 5).Value = '  +2.11%  '
$ws.Cells.Item(42, 5).Value = '  +0.34%  '
$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(43, 4).Formula = "'1.19"
$ws.Cells.Item(43, 5).Value = '  -1.52%  '
$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(44, 4).Formula = "'0.0212"
$ws.Cells.Item(44, 5).Value = '  +0.26%  '
$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(45, 4).Formula = "'16.23"
$ws.Cells.Item(45, 5).Value = '  +3.13%  '
$ws.Cells.Item(46, 2).Value = 'ARBITRUM'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(46, 4).Formula = "'1.04"
$ws.Cells.Item(46, 5).Value = '  -2.79%  '
$ws.Cells.Item(47, 2).Value = 'Maker'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(47, 4).Value = '1.357.87'
$ws.Cells.Item(47, 5).Value = '  +1.94%  '
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(48, 4).Formula = "'88.77"
$ws.Cells.Item(48, 5).Value = '  -2.11%  '
$ws.Cells.Item(49, 2).Value = 'FraxShare'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(49, 4).Formula = "'7.26"
$ws.Cells.Item(49, 5).Value = '  -2.75%  '
$ws.Cells.Item(50, 2).Value = 'MXToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(50, 4).Formula = "'2.82"
$ws.Cells.Item(50, 5).Value = '  -0.13%  '
$ws.Cells.Item(51, 2).Value = 'MultiversX'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(51, 4).Formula = "'46.16"
$ws.Cells.Item(51, 5).Value = '  +5.33%  '
